$d = $word.ActiveDocument

# Locate the paragraph that ends with "Sage is a tool available to you ..."
# so we can insert the new bullet right after it, inheriting its numbering
# (numId=3, ilvl=0) and style.
$anchorText = "Sage is a tool available to you to help you through your time at Cardiff."
$anchorIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $anchorText) {
        $anchorIndex = $i
        $anchorRange = $p.Range
        break
    }
}

if ($anchorIndex -ge 0) {
    # Insert a new paragraph right after the anchor; it copies the anchor's
    # paragraph formatting (including the numPr list numbering).
    $anchorRange.InsertParagraphAfter()

    # The freshly inserted paragraph is now the next one in the collection.
    $newPara = $d.Paragraphs.Item($anchorIndex + 1)
    $newPara.Range.Text = "Sage allows you to share files with particular people (if you know their username) and also allows you to publish it."
}
